$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Black Amber / Primera, date 44210 -> 44230, volume/prices/unit shift to 18 kilos
$ws.Range("D2").Value = 44230
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 16500
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16750
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("S2").Value = 931
$ws.Range("T2").Value = 18

# Row 3: Black Amber / Segunda, date 44210 -> 44230, volume/prices/unit shift to 18 kilos
$ws.Range("D3").Value = 44230
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 14500
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14750
$ws.Range("Q3").Value = "$/caja 18 kilos granel"
$ws.Range("S3").Value = 819
$ws.Range("T3").Value = 18

# Row 4: Black Amber / Primera, date 44230 -> 44209, volume/prices/unit shift to 16 kilos
$ws.Range("D4").Value = 44209
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 15500
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15750
$ws.Range("Q4").Value = "$/caja 16 kilos granel"
$ws.Range("S4").Value = 984
$ws.Range("T4").Value = 16

# Row 5: date 44230 -> 44224, quality Segunda -> Especial, volume/prices/unit shift to 16 kilos
$ws.Range("D5").Value = 44224
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 16500
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 16750
$ws.Range("Q5").Value = "$/caja 16 kilos granel"
$ws.Range("S5").Value = 1047
$ws.Range("T5").Value = 16

# Row 6: quality Especial -> Primera, volume/prices updated
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 14500
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14750
$ws.Range("S6").Value = 922

# Row 7: quality Primera -> Segunda, prices updated
$ws.Range("L7").Value = "Segunda"
$ws.Range("N7").Value = 12500
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12750
$ws.Range("S7").Value = 797

# Row 8: date 44224 -> 44210, quality Segunda -> Primera, volume/prices updated
$ws.Range("D8").Value = 44210
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 15500
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15750
$ws.Range("S8").Value = 984

# Row 9: date 44209 -> 44210, quality Primera -> Segunda, prices updated
$ws.Range("D9").Value = 44210
$ws.Range("L9").Value = "Segunda"
$ws.Range("N9").Value = 12500
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 12750
$ws.Range("S9").Value = 797
